$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix punctuation typos in three "Razon social"/"Nombre Fantasia" entries ---
# (a stray comma was meant to be a period in these names)
$ws.Range("E30").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E39").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E60").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E31").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F31").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E49").Value = "GIMENEZ. ROBERTO ADRIAN"
$ws.Range("F49").Value = "GIMENEZ. ROBERTO ADRIAN"

# --- Fix "Importe" column: these were scraped as locale-formatted text
# ("1.234,56") and must become plain decimal text ("1234.56") --
# the cells stay text, so coerce the column to Text first, write the
# corrected literals, then restore the default (General) style so the
# on-disk cell formatting is unchanged.
$importeRange = $ws.Range("H2:H87")
$importeRange.NumberFormat = "@"
$ws.Range("H2").Value = "25.00"
$ws.Range("H3").Value = "5522.56"
$ws.Range("H4").Value = "292.00"
$ws.Range("H5").Value = "141110.32"
$ws.Range("H6").Value = "1509.50"
$ws.Range("H7").Value = "8150.95"
$ws.Range("H8").Value = "11505.87"
$ws.Range("H9").Value = "1275.00"
$ws.Range("H10").Value = "64.50"
$ws.Range("H11").Value = "14.90"
$ws.Range("H12").Value = "477.50"
$ws.Range("H13").Value = "622.92"
$ws.Range("H14").Value = "81.09"
$ws.Range("H15").Value = "37453.04"
$ws.Range("H16").Value = "33.00"
$ws.Range("H17").Value = "188.00"
$ws.Range("H18").Value = "29.25"
$ws.Range("H19").Value = "1816.62"
$ws.Range("H20").Value = "131162.00"
$ws.Range("H21").Value = "348.32"
$ws.Range("H22").Value = "3560.00"
$ws.Range("H23").Value = "101.74"
$ws.Range("H24").Value = "107.78"
$ws.Range("H25").Value = "1780.17"
$ws.Range("H26").Value = "2000.00"
$ws.Range("H27").Value = "5000.00"
$ws.Range("H28").Value = "750.00"
$ws.Range("H29").Value = "1620.00"
$ws.Range("H30").Value = "1153.00"
$ws.Range("H31").Value = "230.64"
$ws.Range("H32").Value = "6065.00"
$ws.Range("H33").Value = "2480.00"
$ws.Range("H34").Value = "131970.00"
$ws.Range("H35").Value = "166.25"
$ws.Range("H36").Value = "2943.70"
$ws.Range("H37").Value = "50.00"
$ws.Range("H38").Value = "96.00"
$ws.Range("H39").Value = "5936.40"
$ws.Range("H40").Value = "2992.50"
$ws.Range("H41").Value = "34.50"
$ws.Range("H42").Value = "13454.80"
$ws.Range("H43").Value = "1646.00"
$ws.Range("H44").Value = "3800.00"
$ws.Range("H45").Value = "7865.00"
$ws.Range("H46").Value = "250.00"
$ws.Range("H47").Value = "500.00"
$ws.Range("H48").Value = "572.00"
$ws.Range("H49").Value = "400.00"
$ws.Range("H50").Value = "800.00"
$ws.Range("H51").Value = "4000.00"
$ws.Range("H52").Value = "442.01"
$ws.Range("H53").Value = "500.00"
$ws.Range("H54").Value = "950.00"
$ws.Range("H55").Value = "6000.00"
$ws.Range("H56").Value = "400.00"
$ws.Range("H57").Value = "200.00"
$ws.Range("H58").Value = "7095.00"
$ws.Range("H59").Value = "1000.00"
$ws.Range("H60").Value = "3960.00"
$ws.Range("H61").Value = "1790.00"
$ws.Range("H62").Value = "48.56"
$ws.Range("H63").Value = "56.96"
$ws.Range("H64").Value = "97.87"
$ws.Range("H65").Value = "1056.00"
$ws.Range("H66").Value = "2000.24"
$ws.Range("H67").Value = "1712.40"
$ws.Range("H68").Value = "6333.80"
$ws.Range("H69").Value = "9.50"
$ws.Range("H70").Value = "1323.75"
$ws.Range("H71").Value = "865.15"
$ws.Range("H72").Value = "200.00"
$ws.Range("H73").Value = "84800.00"
$ws.Range("H74").Value = "406009.48"
$ws.Range("H75").Value = "201400.00"
$ws.Range("H76").Value = "17000.00"
$ws.Range("H77").Value = "196200.00"
$ws.Range("H78").Value = "10000.00"
$ws.Range("H79").Value = "173961.00"
$ws.Range("H80").Value = "189484.00"
$ws.Range("H81").Value = "207000.00"
$ws.Range("H82").Value = "204524.00"
$ws.Range("H83").Value = "190000.00"
$ws.Range("H84").Value = "444641.06"
$ws.Range("H85").Value = "254255.61"
$ws.Range("H86").Value = "767.65"
$ws.Range("H87").Value = "3700.00"
$importeRange.Style = "Normal"

